$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 84; $r++) {
    $cell = $ws.Range("C$r")
    if ($cell.Value2 -eq 45188) {
        $cell.Value = 45189
    }
}
